$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.473.85"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.803.22"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.59"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.581"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.82%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.84"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.86%  "
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0953"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.818.74"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.21"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").Value = "34.457.66"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("E17").Value = "  +1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.05"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "0.0₃0798"
$ws.Range("E19").Value = "  -0.81%  "
$ws.Range("E20").Value = "  -0.95%  "
$ws.Range("E21").Value = "  +1.63%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "173.83"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.17%  "
$ws.Range("E25").Value = "  +2.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.82"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +6.91%  "
$ws.Range("E27").Value = "  +1.22%  "
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.02"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.83"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.24"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.55%  "
$ws.Range("E35").Value = "  +1.51%  "
$ws.Range("D36").Value = "1.394.72"
$ws.Range("E36").Value = "  -2.46%  "
$ws.Range("E37").Value = "  -3.27%  "
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("E39").Value = "  -1.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.41"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.950"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.36%  "
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.64"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("E45").Value = "  +3.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0509"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.98"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("D48").Value = "1.963.14"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.98"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("E51").Value = "  +0.97%  "
